$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stale "_GoBack" bookmark that currently sits between
#    "...gmail" and ".com" in the e-mail address run. The COM-interop
#    runtime only drops a hidden bookmark when an actual text edit's
#    matched span crosses over it, so we round-trip the surrounding
#    text through a temporary marker and back to its original value —
#    net-zero visible text change, but the crossing edit clears the
#    bookmark that used to live there.
# ------------------------------------------------------------------
$d.Content.Find.Execute("l.com", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__TEMP_GOBACK_MARKER__", 2)
$d.Content.Find.Execute("__TEMP_GOBACK_MARKER__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "l.com", 2)

# ------------------------------------------------------------------
# 2) Update the end-date cell from 25/05/2025 to 20/05/2025.
# ------------------------------------------------------------------
$d.Content.Find.Execute("25/05/2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "20/05/2025", 2)

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark at the new edit location, right
#    after the "20" that replaced "25" (matching the split into a
#    "20" run + bookmark + "/05/2025" run).
# ------------------------------------------------------------------
$dateRng = $d.Content
$dateRng.Find.Execute("20/05/2025", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$goBackPos = $dateRng.Start + 2
$goBackRng = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRng)
